# Auto-generated edit script applying cached-value corrections to the
# per-job-class Leve tables (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR), refreshing
# currentAveragePrice* and Leve*Profit* columns (H:N) to match the latest
# market-board pull from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H34").Value = 1063.4546
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 1063.4546
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H41").Value = 773.9
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H43").Value = 4100.6
$ws.Range("I43").Value = 3333.6667
$ws.Range("K43").Value = 3333.6667
$ws.Range("M43").Value = -3264.6667
$ws.Range("H55").Value = 957.5
$ws.Range("I55").Value = 1122.5
$ws.Range("J55").Value = 875
$ws.Range("K55").Value = 1122.5
$ws.Range("L55").Value = 875
$ws.Range("M55").Value = -908.5
$ws.Range("N55").Value = -1303
$ws.Range("H70").Value = 5599.6665
$ws.Range("I70").Value = 7666.6665
$ws.Range("J70").Value = 4221.6665
$ws.Range("K70").Value = 22999.9995
$ws.Range("L70").Value = 12664.9995
$ws.Range("M70").Value = -22729.9995
$ws.Range("N70").Value = -13204.9995
$ws.Range("H73").Value = 5599.6665
$ws.Range("I73").Value = 7666.6665
$ws.Range("J73").Value = 4221.6665
$ws.Range("K73").Value = 22999.9995
$ws.Range("L73").Value = 12664.9995
$ws.Range("M73").Value = -22063.9995
$ws.Range("N73").Value = -14536.9995
$ws.Range("H74").Value = 5551.143
$ws.Range("I74").Value = 4643
$ws.Range("K74").Value = 4643
$ws.Range("M74").Value = -3707
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("N76").Value = 0
$ws.Range("H77").Value = 5551.143
$ws.Range("I77").Value = 4643
$ws.Range("K77").Value = 23215
$ws.Range("M77").Value = -18535
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("N79").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("N93").Value = 0
$ws.Range("H112").Value = 4966.5
$ws.Range("J112").Value = 4933
$ws.Range("L112").Value = 14799
$ws.Range("N112").Value = -17015
$ws.Range("H116").Value = 7999.375
$ws.Range("I116").Value = 9999
$ws.Range("J116").Value = 7713.7144
$ws.Range("K116").Value = 9999
$ws.Range("L116").Value = 7713.7144
$ws.Range("M116").Value = -6557
$ws.Range("N116").Value = -14597.7144
$ws.Range("H125").Value = 1196.6666
$ws.Range("J125").Value = 1196.6666
$ws.Range("L125").Value = 10769.9994
$ws.Range("N125").Value = -15689.9994

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H11").Value = 177836
$ws.Range("J11").Value = 18000
$ws.Range("L11").Value = 18000
$ws.Range("N11").Value = -18288
$ws.Range("H28").Value = 8051.4287
$ws.Range("I28").Value = 7310
$ws.Range("J28").Value = 12500
$ws.Range("K28").Value = 7310
$ws.Range("L28").Value = 12500
$ws.Range("M28").Value = -7118
$ws.Range("N28").Value = -12884
$ws.Range("H32").Value = 4596.6895
$ws.Range("I32").Value = 4596.6895
$ws.Range("K32").Value = 4596.6895
$ws.Range("M32").Value = -4309.6895
$ws.Range("H99").Value = 8051.4287
$ws.Range("I99").Value = 7310
$ws.Range("J99").Value = 12500
$ws.Range("K99").Value = 7310
$ws.Range("L99").Value = 12500
$ws.Range("M99").Value = -4315
$ws.Range("N99").Value = -18490
$ws.Range("H130").Value = 22222
$ws.Range("J130").Value = 22222
$ws.Range("L130").Value = 22222
$ws.Range("N130").Value = -32262

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H96").Value = 24898
$ws.Range("I96").Value = 24898
$ws.Range("K96").Value = 24898
$ws.Range("M96").Value = -22152

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H137").Value = 40000
$ws.Range("I137").Value = 40000
$ws.Range("K137").Value = 40000
$ws.Range("M137").Value = -34900

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 3876700
$ws.Range("I4").Value = 1804529
$ws.Range("K4").Value = 5413587
$ws.Range("M4").Value = -5413475
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 1800
$ws.Range("M22").Value = -1031
$ws.Range("N22").Value = -2138
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = 0
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 600
$ws.Range("K27").Value = 1200
$ws.Range("L27").Value = 1800
$ws.Range("M27").Value = -1098
$ws.Range("N27").Value = -2004
$ws.Range("H86").Value = 7999
$ws.Range("J86").Value = 7999
$ws.Range("L86").Value = 23997
$ws.Range("N86").Value = -26369
$ws.Range("H88").Value = 20000
$ws.Range("J88").Value = 20000
$ws.Range("L88").Value = 60000
$ws.Range("N88").Value = -60856
$ws.Range("H89").Value = 7999
$ws.Range("J89").Value = 7999
$ws.Range("L89").Value = 71991
$ws.Range("N89").Value = -83847
$ws.Range("H91").Value = 20000
$ws.Range("J91").Value = 20000
$ws.Range("L91").Value = 60000
$ws.Range("N91").Value = -62964
$ws.Range("H99").Value = 1084.5
$ws.Range("I99").Value = 1436.6666
$ws.Range("J99").Value = 28
$ws.Range("K99").Value = 4309.9998
$ws.Range("L99").Value = 84
$ws.Range("M99").Value = -2063.9998
$ws.Range("N99").Value = -4576
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H107").Value = 137.25
$ws.Range("I107").Value = 97.5
$ws.Range("J107").Value = 177
$ws.Range("K107").Value = 292.5
$ws.Range("L107").Value = 531
$ws.Range("M107").Value = 1627.5
$ws.Range("N107").Value = -4371
$ws.Range("H108").Value = 487.5
$ws.Range("I108").Value = 487.5
$ws.Range("K108").Value = 1462.5
$ws.Range("M108").Value = 1417.5
$ws.Range("H119").Value = 1599.2
$ws.Range("I119").Value = 1499
$ws.Range("J119").Value = 2000
$ws.Range("K119").Value = 4497
$ws.Range("L119").Value = 6000
$ws.Range("M119").Value = 341
$ws.Range("N119").Value = -15676

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 7861
$ws.Range("I70").Value = 7224.5
$ws.Range("J70").Value = 8497.5
$ws.Range("K70").Value = 7224.5
$ws.Range("L70").Value = 8497.5
$ws.Range("M70").Value = -6954.5
$ws.Range("N70").Value = -9037.5
$ws.Range("H73").Value = 7861
$ws.Range("I73").Value = 7224.5
$ws.Range("J73").Value = 8497.5
$ws.Range("K73").Value = 7224.5
$ws.Range("L73").Value = 8497.5
$ws.Range("M73").Value = -6288.5
$ws.Range("N73").Value = -10369.5
$ws.Range("H122").Value = 4757.7
$ws.Range("I122").Value = 3666.3333
$ws.Range("J122").Value = 5225.4287
$ws.Range("K122").Value = 10998.9999
$ws.Range("L122").Value = 15676.2861
$ws.Range("M122").Value = -8548.999899999999
$ws.Range("N122").Value = -20576.2861
$ws.Range("H132").Value = 3561.8572
$ws.Range("I132").Value = 3503.1667
$ws.Range("K132").Value = 10509.5001
$ws.Range("M132").Value = -7979.500100000001

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H55").Value = 833
$ws.Range("I55").Value = 405
$ws.Range("K55").Value = 405
$ws.Range("M55").Value = -232
$ws.Range("H98").Value = 48529.89
$ws.Range("J98").Value = 48529.89
$ws.Range("L98").Value = 48529.89
$ws.Range("N98").Value = -54519.89
$ws.Range("H122").Value = 3553.7407
$ws.Range("I122").Value = 3798.4
$ws.Range("K122").Value = 11395.2
$ws.Range("M122").Value = -8945.200000000001

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H96").Value = 5138.8
$ws.Range("I96").Value = 4898.8335
$ws.Range("K96").Value = 4898.8335
$ws.Range("M96").Value = -3525.8335
$ws.Range("H122").Value = 2676.818
$ws.Range("I122").Value = 2637.7778
$ws.Range("J122").Value = 2852.5
$ws.Range("K122").Value = 7913.3334
$ws.Range("L122").Value = 8557.5
$ws.Range("M122").Value = -5463.3334
$ws.Range("N122").Value = -13457.5
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").ClearContents()
$ws.Range("N129").Value = 0
